$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows data: (RowNum, A, B, C, D(date serial), E, F, G, H, I, J, K, L, M, N, O, P, Q, R)
$rows = @(
    @{ RowNum=70; A=1; B="Agrícola del Norte S.A. de Arica"; C="Arica y Parinacota"; D=44911; E=15; F=100112027; G="Melón"; H="Tuna"; I="Primera";  J=100; K=19000; L=20000; M=19500; N="$/caja 18 unidades"; O="Región de Arica y Parinacota"; P=1083; Q=18; R="Hortaliza" },
    @{ RowNum=71; A=1; B="Agrícola del Norte S.A. de Arica"; C="Arica y Parinacota"; D=44911; E=15; F=100112027; G="Melón"; H="Tuna"; I="Segunda"; J=70;  K=16000; L=17000; M=16429; N="$/caja 24 unidades"; O="Región de Arica y Parinacota"; P=685;  Q=24; R="Hortaliza" },
    @{ RowNum=72; A=1; B="Agrícola del Norte S.A. de Arica"; C="Arica y Parinacota"; D=44911; E=15; F=100112027; G="Melón"; H="Tuna"; I="Tercera";  J=60;  K=14000; L=15000; M=14667; N="$/caja 30 unidades"; O="Región de Arica y Parinacota"; P=489;  Q=30; R="Hortaliza" }
)

foreach ($row in $rows) {
    $r = $row.RowNum
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
}
